$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial 45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update prices in D29 and D30
$ws.Range("D29").Value = 520.458
$ws.Range("D30").Value = 353.073
